$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet/tab name to reflect new "through" date
$ws.Name = "Through 2021-12-13"

# Row 12 (October) updates
$ws.Range("U12").Value = 188
$ws.Range("V12").Value = 0.0408

# Row 14 (December) updates
$ws.Range("A14").Value = "December (through 12-13)"
$ws.Range("C14").Value = 11
$ws.Range("D14").Value = 0.1538
$ws.Range("F14").Value = 34
$ws.Range("G14").Value = 0.1053
$ws.Range("I14").Value = 41
$ws.Range("J14").Value = 0.0682
$ws.Range("L14").Value = 25
$ws.Range("M14").Value = 0.1071
$ws.Range("O14").Value = 19
$ws.Range("P14").Value = 0.1364
$ws.Range("R14").Value = 63
$ws.Range("S14").Value = 0.0455
$ws.Range("U14").Value = 95

# Row 15 (Total) updates
$ws.Range("C15").Value = 269
$ws.Range("D15").Value = 0.1151
$ws.Range("F15").Value = 537
$ws.Range("G15").Value = 0.1065
$ws.Range("I15").Value = 799
$ws.Range("J15").Value = 0.0763
$ws.Range("L15").Value = 633
$ws.Range("M15").Value = 0.1085
$ws.Range("O15").Value = 499
$ws.Range("P15").Value = 0.1025
$ws.Range("R15").Value = 1263
$ws.Range("S15").Value = 0.0504
$ws.Range("U15").Value = 1638
$ws.Range("V15").Value = 0.0575
